$d = $word.ActiveDocument

$oldText = "At two given times two of the passengers (cat, parrot, or seed) are on one side of the river leaving one passenger on the other side of the river. The man can only be in one place at a given time."
$newText = "Now the sub goals are the delivery of each parcel individually while leaving two parcels unprotected at two given times."

# Locate the paragraph containing the sentence that needs to be split into two
# paragraphs (with a blank paragraph between them).
$findRange = $d.Content
$found = $findRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target sentence in document."
}

$paraRange = $findRange.Duplicate
$paraRange.Expand(4) | Out-Null

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$newXml = '<w:p xmlns:w="' + $w + '">' +
            '<w:r><w:tab/></w:r>' +
            '<w:r><w:t>' + $oldText + '</w:t></w:r>' +
          '</w:p>' +
          '<w:p xmlns:w="' + $w + '"/>' +
          '<w:p xmlns:w="' + $w + '">' +
            '<w:r><w:tab/></w:r>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
            '<w:bookmarkEnd w:id="0"/>' +
            '<w:r><w:t>' + $newText + '</w:t></w:r>' +
          '</w:p>'

$paraRange.InsertXML($newXml)
